# Update the per-seed patient scores (and per-seed/overall averages) in the
# "seed_results_CTNet_balanced" results sheet so that they are reported with
# two decimal places of precision instead of three, for columns
# Patient1..Patient9 and Average (B:K), rows 2-16 (seed rows 2-15, Average row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> values for columns B,C,D,E,F,G,H,I,J,K
$data = @{
    2  = @(0.77, 0.71, 0.76, 0.97, 0.91, 0.87, 0.91, 0.94, 0.88, 0.86)
    3  = @(0.81, 0.70, 0.77, 0.90, 0.98, 0.86, 0.91, 0.94, 0.90, 0.86)
    4  = @(0.79, 0.63, 0.78, 0.88, 0.97, 0.93, 0.91, 0.95, 0.89, 0.86)
    5  = @(0.79, 0.71, 0.77, 0.90, 0.98, 0.88, 0.90, 0.94, 0.88, 0.86)
    6  = @(0.75, 0.69, 0.74, 0.86, 0.97, 0.87, 0.92, 0.94, 0.89, 0.85)
    7  = @(0.76, 0.69, 0.79, 0.93, 0.99, 0.87, 0.90, 0.95, 0.91, 0.87)
    8  = @(0.81, 0.71, 0.80, 0.93, 0.97, 0.86, 0.90, 0.94, 0.90, 0.87)
    9  = @(0.77, 0.72, 0.72, 0.96, 0.97, 0.88, 0.91, 0.95, 0.89, 0.86)
    10 = @(0.77, 0.69, 0.78, 0.90, 0.98, 0.87, 0.92, 0.93, 0.88, 0.86)
    11 = @(0.76, 0.67, 0.77, 0.97, 0.97, 0.85, 0.91, 0.95, 0.90, 0.86)
    12 = @(0.77, 0.70, 0.71, 0.90, 0.98, 0.87, 0.89, 0.94, 0.90, 0.85)
    13 = @(0.73, 0.69, 0.75, 0.97, 0.96, 0.90, 0.89, 0.95, 0.87, 0.86)
    14 = @(0.76, 0.62, 0.76, 0.91, 0.96, 0.89, 0.91, 0.94, 0.89, 0.85)
    15 = @(0.80, 0.69, 0.75, 0.96, 0.96, 0.90, 0.92, 0.94, 0.88, 0.87)
    16 = @(0.77, 0.69, 0.76, 0.92, 0.97, 0.88, 0.91, 0.94, 0.89, 0.86)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # column B is the 2nd column
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
